# edit.ps1 -- reproduces the "added calc of max working hour per month ..." commit
# against invoicer/data_result/result.xlsx (sheet "WorkTimeReport").
#
# Summary of the edit:
#  - the weekend-styled date cells (style index 1, green fill) are bumped
#    from April to May (only those rows -- the normal working-day rows keep
#    their original April dates, matching the source diff exactly)
#  - most "Task description" (column F) cells get more descriptive text
#  - a new weekend row (31.05.2024) is appended as row 32
#  - the "total" row (SUM of column D) moves down from row 32 to row 33 and
#    now sums D1:D32 instead of D1:D31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a literal date-look-alike string (dd.mm.yyyy) into a cell
# without Excel's autodetection turning it into a real date serial number.
# We flip the cell to Text format just long enough to stuff the literal
# string in, matching what a human pasting "Text"-formatted data would do.
# ---------------------------------------------------------------------
function Set-DateText($a1, [string]$text) {
    $cell = $ws.Range($a1)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

# ---- weekend-row dates: April -> May -------------------------------------
Set-DateText "A2"  "01.05.2024"
Set-DateText "A7"  "06.05.2024"
Set-DateText "A8"  "07.05.2024"
Set-DateText "A14" "13.05.2024"
Set-DateText "A15" "14.05.2024"
Set-DateText "A21" "20.05.2024"
Set-DateText "A22" "21.05.2024"
Set-DateText "A28" "27.05.2024"
Set-DateText "A29" "28.05.2024"

# ---- task description (column F) rewordings -------------------------------
$ws.Range("F3").Value2  = "test automation for admin part`n"
$ws.Range("F4").Value2  = "pi planning day 1`n"
$ws.Range("F5").Value2  = "pi planning day 2`n"
$ws.Range("F6").Value2  = "creating pipeline for regression tests`n"
$ws.Range("F9").Value2  = "test automations for panels refactor`n"
$ws.Range("F10").Value2 = "tests for admin acquirer and messages`n"
$ws.Range("F11").Value2 = "refactor for waits`n"
$ws.Range("F12").Value2 = "test automation refactor for waits and unused methods`n"
$ws.Range("F13").Value2 = "retest of bug for incorrect transaction failed message`n"
$ws.Range("F16").Value2 = "tests refactor`n"
$ws.Range("F17").Value2 = "prod deploy sanity tests, tests automation`n"
$ws.Range("F18").Value2 = "automations for portal tests`n"
$ws.Range("F19").Value2 = "refactor of automated tests`n"
$ws.Range("F20").Value2 = "tests smell code fixes, updating repo`n"
$ws.Range("F23").Value2 = "refactor tests for class methods`n"
$ws.Range("F24").Value2 = "test automation for transaction details`n"
$ws.Range("F25").Value2 = "transaction advance filters automated tests, pipeline clarification meeting`n"
$ws.Range("F26").Value2 = "reseeding db meeting, automation tests`n"
$ws.Range("F27").Value2 = "investigation for pipeline bugs, automation tests`n"
$ws.Range("F30").Value2 = "terminal tests refactor`n"
$ws.Range("F31").Value2 = "admin tests refactor`n"

# ---- the old total row (row 32, D32=SUM(D1:D31)) is replaced by a new
#      weekend data row, and the total moves to row 33 --------------------
$ws.Range("D32").ClearContents()

# give the new row the same "weekend" look (green fill) as the other
# weekend rows, e.g. A7:F7, *before* filling in its values so the new
# text-formatted date cell picks up the fill too
$ws.Range("A7:F7").Copy()
$ws.Range("A32:F32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

Set-DateText "A32" "31.05.2024"
$ws.Range("B32").Value2 = "John Doe"
$ws.Range("C32").Value2 = "TestClient"
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = 0

$ws.Range("D33").Formula = "=SUM(D1:D32)"
